# Generate Report for Handback
#
# Refreshing the handback report updates the handoff/handback generation
# timestamps recorded for the file "8923af4e-e3ce-4157-bf36-7321614f208a.md"
# (row 4 in each sheet) across the Overview sheet and each locale sheet
# (zh-cn, de-de).

$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" column (G) for this file
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G4").Value = "2016-12-16 08:42:53"

# zh-cn sheet: Correspond Handoff Datetime (H) / Correspond Handback DateTime (L)
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H4").Value = "2016-12-16 08:42:39"
$zhcn.Range("L4").Value = "2016-12-16 08:43:32"

# de-de sheet: Correspond Handoff Datetime (H) / Correspond Handback DateTime (L)
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H4").Value = "2016-12-16 08:42:53"
$dede.Range("L4").Value = "2016-12-16 08:43:50"
